$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 26-49: refreshed stimuli list (image, trialnum, condition, word, location, repetition)
$rows = @(
    @("E.png", 25, "S", 4, 1, 2),
    @("L.png", 26, "S", 4, 2, 2),
    @("H.png", 27, "S", 4, 3, 2),
    @("A.png", 28, "S", 2, 1, 3),
    @("J.png", 29, "S", 2, 2, 3),
    @("F.png", 30, "S", 2, 3, 3),
    @("E.png", 31, "S", 4, 1, 3),
    @("L.png", 32, "S", 4, 2, 3),
    @("H.png", 33, "S", 4, 3, 3),
    @("D.png", 34, "S", 3, 1, 3),
    @("G.png", 35, "S", 3, 2, 3),
    @("C.png", 36, "S", 3, 3, 3),
    @("M.png", 37, "S", 1, 1, 4),
    @("K.png", 38, "S", 1, 2, 4),
    @("B.png", 39, "S", 1, 3, 4),
    @("D.png", 40, "S", 3, 1, 4),
    @("G.png", 41, "S", 3, 2, 4),
    @("C.png", 42, "S", 3, 3, 4),
    @("M.png", 43, "S", 1, 1, 5),
    @("K.png", 44, "S", 1, 2, 5),
    @("B.png", 45, "S", 1, 3, 5),
    @("E.png", 46, "S", 4, 1, 4),
    @("L.png", 47, "S", 4, 2, 4),
    @("H.png", 48, "S", 4, 3, 4)
)

$r = 26
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Restore the active selection to match the refreshed data extent
[void]$ws.Range("A26:F49").Select()
